$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Percent-formatted humidity cells need an explicit Text number format first,
# otherwise Excel auto-converts "NN%" strings into a numeric percentage value.
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H17").NumberFormat = "@"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H27").NumberFormat = "@"
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H42").NumberFormat = "@"
$ws.Range("H45").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-17 18:18:21"
$ws.Range("O2").Value = "2.6 °C"
$ws.Range("E3").Value = "2026-02-17 18:18:23"
$ws.Range("O3").Value = "-4.3 °C"
$ws.Range("E4").Value = "2026-02-17 18:18:25"
$ws.Range("H4").Value = "81%"
$ws.Range("J4").Value = "1018.2 hPa"
$ws.Range("E5").Value = "2026-02-17 18:18:27"
$ws.Range("I5").Value = "2.0 mm"
$ws.Range("O5").Value = "-3.8 °C"
$ws.Range("E6").Value = "2026-02-17 18:18:30"
$ws.Range("J6").Value = "1018.2 hPa"
$ws.Range("O6").Value = "10.7 °C"
$ws.Range("E7").Value = "2026-02-17 18:18:32"
$ws.Range("H7").Value = "59%"
$ws.Range("E8").Value = "2026-02-17 18:18:34"
$ws.Range("H8").Value = "71%"
$ws.Range("J8").Value = "1018.0 hPa"
$ws.Range("E9").Value = "2026-02-17 18:18:37"
$ws.Range("E10").Value = "2026-02-17 18:18:39"
$ws.Range("H10").Value = "72%"
$ws.Range("O10").Value = "10.8 °C"
$ws.Range("E11").Value = "2026-02-17 18:18:41"
$ws.Range("H11").Value = "47%"
$ws.Range("E12").Value = "2026-02-17 18:18:44"
$ws.Range("H12").Value = "56%"
$ws.Range("E13").Value = "2026-02-17 18:18:46"
$ws.Range("J13").Value = "1017.6 hPa"
$ws.Range("E14").Value = "2026-02-17 18:18:48"
$ws.Range("H14").Value = "65%"
$ws.Range("O14").Value = "13.6 °C"
$ws.Range("E15").Value = "2026-02-17 18:18:51"
$ws.Range("H15").Value = "55%"
$ws.Range("N15").Value = "9.6 °C 17:59 TU"
$ws.Range("E16").Value = "2026-02-17 18:18:53"
$ws.Range("K16").Value = "9.8 MJ/m2"
$ws.Range("E17").Value = "2026-02-17 18:18:55"
$ws.Range("H17").Value = "79%"
$ws.Range("E18").Value = "2026-02-17 18:18:58"
$ws.Range("H18").Value = "78%"
$ws.Range("J18").Value = "1018.4 hPa"
$ws.Range("O18").Value = "10.4 °C"
$ws.Range("E19").Value = "2026-02-17 18:19:00"
$ws.Range("H19").Value = "72%"
$ws.Range("E20").Value = "2026-02-17 18:19:02"
$ws.Range("H20").Value = "60%"
$ws.Range("E21").Value = "2026-02-17 18:19:04"
$ws.Range("O21").Value = "9.9 °C"
$ws.Range("E22").Value = "2026-02-17 18:19:07"
$ws.Range("E23").Value = "2026-02-17 18:19:09"
$ws.Range("I23").Value = "2.5 mm"
$ws.Range("O23").Value = "-3.9 °C"
$ws.Range("E24").Value = "2026-02-17 18:19:12"
$ws.Range("E25").Value = "2026-02-17 18:19:14"
$ws.Range("H25").Value = "50%"
$ws.Range("O25").Value = "-0.9 °C"
$ws.Range("E26").Value = "2026-02-17 18:19:16"
$ws.Range("E27").Value = "2026-02-17 18:19:18"
$ws.Range("H27").Value = "51%"
$ws.Range("O27").Value = "-0.5 °C"
$ws.Range("E28").Value = "2026-02-17 18:19:21"
$ws.Range("J28").Value = "1018.0 hPa"
$ws.Range("E29").Value = "2026-02-17 18:19:23"
$ws.Range("O29").Value = "12.4 °C"
$ws.Range("E30").Value = "2026-02-17 18:19:26"
$ws.Range("H30").Value = "60%"
$ws.Range("J30").Value = "1018.1 hPa"
$ws.Range("E31").Value = "2026-02-17 18:19:28"
$ws.Range("E32").Value = "2026-02-17 18:19:30"
$ws.Range("E33").Value = "2026-02-17 18:19:32"
$ws.Range("H33").Value = "40%"
$ws.Range("E34").Value = "2026-02-17 18:19:35"
$ws.Range("E35").Value = "2026-02-17 18:19:37"
$ws.Range("H35").Value = "92%"
$ws.Range("J35").Value = "1019.9 hPa"
$ws.Range("E36").Value = "2026-02-17 18:19:40"
$ws.Range("H36").Value = "57%"
$ws.Range("J36").Value = "1018.3 hPa"
$ws.Range("O36").Value = "12.6 °C"
$ws.Range("E37").Value = "2026-02-17 18:19:42"
$ws.Range("J37").Value = "1018.7 hPa"
$ws.Range("E38").Value = "2026-02-17 18:19:44"
$ws.Range("O38").Value = "11.3 °C"
$ws.Range("E39").Value = "2026-02-17 18:19:47"
$ws.Range("H39").Value = "56%"
$ws.Range("E40").Value = "2026-02-17 18:19:49"
$ws.Range("H40").Value = "50%"
$ws.Range("O40").Value = "9.7 °C"
$ws.Range("E41").Value = "2026-02-17 18:19:51"
$ws.Range("O41").Value = "16.9 °C"
$ws.Range("E42").Value = "2026-02-17 18:19:54"
$ws.Range("H42").Value = "56%"
$ws.Range("E43").Value = "2026-02-17 18:19:56"
$ws.Range("E44").Value = "2026-02-17 18:19:58"
$ws.Range("O44").Value = "-3.2 °C"
$ws.Range("E45").Value = "2026-02-17 18:20:01"
$ws.Range("H45").Value = "66%"
$ws.Range("I45").Value = "0.3 mm"
$ws.Range("E46").Value = "2026-02-17 18:20:03"
